{"js": "// The underlying change recorded in the template's canonical OOXML is a\n// pure re-serialization: every element in word/document.xml (the root\n// namespace declarations, <w:pgSz>, <w:pgMar>) and in word/styles.xml\n// (<w:rFonts>, <w:lang>, <w:latentStyles>, every <w:lsdException>, the\n// <w:style> definitions for Normal / Default Paragraph Font / Normal\n// Table / No List, <w:tblInd> and <w:tblCellMar>) simply had its XML\n// attributes re-ordered (alphabetically) by whatever tool re-saved the\n// resource \u2014 not a single attribute value, piece of text, or style\n// setting actually changed. (Diffing the canonicalized/C14N form of the\n// before and after XML confirms they are byte-for-byte identical.)\n//\n// Word's JS API has no \"resort the raw attributes\" primitive (and\n// shouldn't), so we reproduce the author's intent the only meaningful\n// way the object model allows: by touching every style that the diff\n// rewrites, re-asserting the exact same values it already has. This\n// exercises the same four <w:style> definitions the diff rewrites\n// without altering anything a reader (or Word) would observe.\n\nconst styles = context.document.getStyles();\nstyles.load(\"items/nameLocal,items/quickStyle,items/unhideWhenUsed\");\nawait context.sync();\n\nfor (const style of styles.items) {\n  // \"Normal\" only carries <w:qFormat/> (no <w:unhideWhenUsed/> in the\n  // XML), so round-trip that flag for it; the other three built-in\n  // styles referenced by the diff (Default Paragraph Font / Normal\n  // Table / No List) already carry an explicit <w:unhideWhenUsed/>, so\n  // round-trip that flag for them instead. Writing back the identical\n  // value touches the style definition without introducing (or\n  // dropping) any attribute.\n  if (style.nameLocal === \"Normal\") {\n    style.quickStyle = style.quickStyle;\n  } else {\n    style.unhideWhenUsed = style.unhideWhenUsed;\n  }\n}\n\nawait context.sync();\n", "ps1": "# The underlying change recorded in the template's canonical OOXML is a\n# pure re-serialization: every element in word/document.xml (the root\n# namespace declarations, <w:pgSz>, <w:pgMar>) and in word/styles.xml\n# (<w:rFonts>, <w:lang>, <w:latentStyles>, every <w:lsdException>, the\n# <w:style> definitions for Normal / Default Paragraph Font / Normal\n# Table / No List, <w:tblInd> and <w:tblCellMar>) simply had its XML\n# attributes re-ordered (alphabetically) by whatever tool re-saved the\n# resource - not a single attribute value, piece of text, or style\n# setting actually changed. (Diffing the canonicalized/C14N form of the\n# before and after XML confirms they are byte-for-byte identical.)\n#\n# The Word object model has no \"resort the raw attributes\" primitive\n# (and shouldn't), so we reproduce the author's intent the only\n# meaningful way it allows: by touching every style the diff rewrites,\n# re-asserting the exact same name it already has. This exercises the\n# same four <w:style> definitions the diff rewrites without altering\n# anything a reader (or Word) would observe.\n\n$d = $word.ActiveDocument\n\nforeach ($name in @(\"Normal\", \"Default Paragraph Font\", \"Normal Table\", \"No List\")) {\n    $s = $d.Styles($name)\n    $s.NameLocal = $s.NameLocal\n}\n"}
